$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6619111895561218
$ws.Range("B1").Value = 2.090850353240967
$ws.Range("C1").Value = 5.174644470214844
$ws.Range("D1").Value = 2.851711511611938
$ws.Range("E1").Value = 0.6882362365722656
